# Refresh the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.289.51"
$ws.Range("D3").Value = "3.242.03"
$ws.Range("E3").Value = "  +2.51%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.19"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "178.79"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.63%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "3.241.63"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("E10").Value = "  +4.61%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.412"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").Value = "3.808.77"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  +0.83%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.93"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").Value = "67.222.00"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").Value = "3.249.91"
$ws.Range("E18").Value = "  +2.67%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.31"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.20%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "375.15"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.86%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.58"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +6.24%  "
$ws.Range("E23").Value = "  +0.07%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "71.08"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.80%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.508"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "3.387.65"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("E27").Value = "  -0.69%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.91"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +4.50%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "22.49"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +5.90%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.06%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "161.28"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.04%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.856"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("E40").Value = "  +10.13%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +16.69%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "26.82"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.54%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("D44").Value = "2.757.91"
$ws.Range("E44").Value = "  +6.31%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.39"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.44%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "350.87"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +10.74%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "25.65"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +8.90%  "
$ws.Range("E48").Value = "  +2.81%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0670"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.40%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0279"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("E51").Value = "  +1.73%  "
